# Actualización automática 2025-08-05 10:15:08
# Updates the "PRESUPUESTO" (column G) values on the "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G3").Value = 1900
$ws.Range("G4").Value = 15100
$ws.Range("G5").Value = 9000
$ws.Range("G6").Value = 2450
$ws.Range("G7").Value = 2450
$ws.Range("G8").Value = 750
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 1350
$ws.Range("G11").Value = 1300
$ws.Range("G13").Value = 650
$ws.Range("G14").Value = 5100
$ws.Range("G15").Value = 3935.11
$ws.Range("G16").Value = 7410
$ws.Range("G18").Value = 4330
$ws.Range("G20").Value = 950
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 1500
$ws.Range("G24").Value = 300
$ws.Range("G25").Value = 16110
$ws.Range("G26").Value = 300
$ws.Range("G29").Value = 5850
$ws.Range("G30").Value = 260
$ws.Range("G31").Value = 350
$ws.Range("G32").Value = 3100
$ws.Range("G34").Value = 0
$ws.Range("G36").Value = 18200
$ws.Range("G37").Value = 5600
$ws.Range("G38").Value = 100
$ws.Range("G41").Value = 5505
$ws.Range("G42").Value = 1020
$ws.Range("G43").Value = 2250
$ws.Range("G44").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("G47").Value = 3300
$ws.Range("G48").Value = 3750
$ws.Range("G49").Value = 300
$ws.Range("G50").Value = 550
$ws.Range("G53").Value = 1650
$ws.Range("G54").Value = 1000
$ws.Range("G55").Value = 0
$ws.Range("G56").Value = 130170.11
